# Update "想去人数" (want-to-go count) figures that changed between scrapes.
# Affects the "展览" sheet and the "全部类型" sheet (the latter aggregates
# all event types, including the same rows as "展览").

$wb = $excel.ActiveWorkbook

$sheetExhibit = $wb.Worksheets.Item("展览")
$sheetExhibit.Range("F2").Value = 7066
$sheetExhibit.Range("F5").Value = 114
$sheetExhibit.Range("F6").Value = 1089
$sheetExhibit.Range("F7").Value = 174
$sheetExhibit.Range("F8").Value = 31

$sheetAll = $wb.Worksheets.Item("全部类型")
$sheetAll.Range("F2").Value = 7066
$sheetAll.Range("F5").Value = 114
$sheetAll.Range("F6").Value = 1089
$sheetAll.Range("F7").Value = 174
$sheetAll.Range("F9").Value = 31
